$p = $ppt.ActivePresentation

# The table on slide 16 had its table style changed (PowerPoint's
# "Table Design" gallery) from the deck's custom "Table_0" style to a
# built-in table style.
$s = $p.Slides.Item(16)

for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $shape = $s.Shapes.Item($i)
    if ($shape.HasTable) {
        $shape.Table.ApplyStyle("{AC7B3A33-4158-4E7A-8461-A2D4A8E7503E}")
    }
}
